$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.914.78'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.774.88'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.89'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.95'
$ws.Range("E8").Value = '  +3.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.288'
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0688'
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0937'
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("D12").Value = '2.033.37'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.92'
$ws.Range("E13").Value = '  +4.60%  '
$ws.Range("D14").Value = '1.765.69'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '33.920.81'
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.614'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.14'
$ws.Range("E17").Value = '  -0.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.18'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.50'
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").Value = '0.0₃0775'
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.62'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.07'
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("E24").Value = '  -2.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.19'
$ws.Range("E25").Value = '  +1.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.15'
$ws.Range("E26").Value = '  -0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.04'
$ws.Range("E27").Value = '  +1.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.112'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  +3.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0509'
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.61'
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.53'
$ws.Range("E33").Value = '  +1.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.78'
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("D35").Value = '1.383.27'
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.654'
$ws.Range("E36").Value = '  +3.65%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0185'
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.22'
$ws.Range("E39").Value = '  +6.01%  '
$ws.Range("E40").Value = '  +0.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.905'
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.65'
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '77.31'
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("E44").Value = '  +22.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.30'
$ws.Range("E45").Value = '  +13.01%  '
$ws.Range("E46").Value = '  +3.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '108.01'
$ws.Range("E47").Value = '  +4.19%  '
$ws.Range("E48").Value = '  +2.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.82'
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("D50").Value = '1.932.95'
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("E51").Value = '  +0.53%  '
